$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Formula = "=""" + $text + """"
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

Set-TextValue "D2" "287.51"
Set-TextValue "E2" "1.57%"
Set-TextValue "D3" "29.51"
Set-TextValue "E3" "4.13%"
Set-TextValue "E4" "0.88%"
Set-TextValue "D5" "0.06752"
Set-TextValue "E5" "3.63%"
Set-TextValue "D6" "7.366"
Set-TextValue "E6" "1.88%"
Set-TextValue "D7" "1.383"
Set-TextValue "E7" "0.02%"
Set-TextValue "D8" "0.9155"
Set-TextValue "E8" "-0.32%"
Set-TextValue "D9" "0.1589"
Set-TextValue "E9" "3.39%"
Set-TextValue "D10" "0.06944"
Set-TextValue "E10" "5.33%"
Set-TextValue "D11" "0.07658"
Set-TextValue "E11" "1.22%"
Set-TextValue "D12" "0.02926"
Set-TextValue "E12" "5.37%"
Set-TextValue "D13" "0.08987"
Set-TextValue "E13" "0.16%"
Set-TextValue "D14" "0.001587"
Set-TextValue "E14" "-0.21%"
Set-TextValue "D15" "0.04481"
Set-TextValue "E15" "1.15%"
Set-TextValue "D16" "0.0006454"
Set-TextValue "E16" "1.20%"
Set-TextValue "D17" "0.006137"
Set-TextValue "E17" "-0.57%"
Set-TextValue "D18" "3.453"
Set-TextValue "E18" "0.18%"
Set-TextValue "D19" "3.445"
Set-TextValue "E20" "-0.21%"
Set-TextValue "E21" "0.50%"
Set-TextValue "D22" "0.1314"
Set-TextValue "E22" "2.53%"
Set-TextValue "D23" "4.099"
Set-TextValue "E23" "2.82%"
Set-TextValue "E24" "2.68%"
Set-TextValue "D25" "0.001196"
Set-TextValue "E25" "1.23%"
Set-TextValue "D26" "0.004139"
Set-TextValue "E26" "-6.72%"
Set-TextValue "E27" "0.13%"
Set-TextValue "D28" "0.0001617"
Set-TextValue "E28" "-0.01%"
Set-TextValue "D40" "0.04256"
Set-TextValue "E40" "3.36%"
Set-TextValue "D41" "0.006819"
Set-TextValue "E41" "2.18%"
Set-TextValue "D42" "0.1240"
Set-TextValue "E42" "0.68%"
Set-TextValue "D43" "0.002232"
Set-TextValue "E43" "3.85%"
Set-TextValue "D44" "0.01292"
Set-TextValue "E44" "6.73%"
Set-TextValue "D45" "0.00005738"
Set-TextValue "E45" "1.17%"
Set-TextValue "D46" "1.960"
Set-TextValue "E46" "-0.30%"
